$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "include?" column (C) gets marked "y" for these MSO rows (AT&T, Charter,
# Comcast, Cox, Verizon) -- fixes the pd.duplicated bug that had been
# skipping Charter et al.
$ws.Range("C2").Value = "y"
$ws.Range("C4").Value = "y"
$ws.Range("C5").Value = "y"
$ws.Range("C6").Value = "y"
$ws.Range("C7").Value = "y"

# C2 previously used a bold font (left over from older formatting); drop the
# bold so it matches the plain centered style already used by the other
# "include?" cells.
$ws.Range("C2").Font.Bold = $false

# Move the active selection to A7.
[void]$ws.Range("A7").Select()
